# Apply cryptos list update (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue 'D2' '29.272.62'
Set-TextValue 'E2' '  +0.50%  '
Set-TextValue 'D3' '1.860.69'
Set-TextValue 'E4' '  +0.06%  '
Set-TextValue 'D5' '0.7037'
Set-TextValue 'E5' '  -0.52%  '
Set-TextValue 'D6' '237.75'
Set-TextValue 'E6' '  -0.22%  '
Set-TextValue 'D8' '0.08210'
Set-TextValue 'E8' '  +9.15%  '
Set-TextValue 'E9' '  -0.36%  '
Set-TextValue 'D10' '23.29'
Set-TextValue 'E10' '  -0.47%  '
Set-TextValue 'D11' '0.08166'
Set-TextValue 'E11' '  +0.40%  '
Set-TextValue 'D12' '1.873.21'
Set-TextValue 'E12' '  -0.61%  '
Set-TextValue 'D13' '0.7164'
Set-TextValue 'E13' '  -1.26%  '
Set-TextValue 'D14' '5.174'
Set-TextValue 'E14' '  -0.91%  '
Set-TextValue 'D15' '89.28'
Set-TextValue 'E15' '  +0.06%  '
Set-TextValue 'D16' '29.291.89'
Set-TextValue 'E16' '  +0.06%  '
Set-TextValue 'D17' '5.778'
Set-TextValue 'E17' '  -0.33%  '
Set-TextValue 'B18' 'Avalanche'
Set-TextValue 'C18' 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue 'D18' '13.39'
Set-TextValue 'E18' '  +2.26%  '
Set-TextValue 'B19' 'ShibaInu'
Set-TextValue 'C19' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue 'D19' '0.000007854'
Set-TextValue 'E19' '  +2.26%  '
Set-TextValue 'D20' '237.57'
Set-TextValue 'E20' '  -0.79%  '
Set-TextValue 'D21' '0.9995'
Set-TextValue 'E21' '  -0.14%  '
Set-TextValue 'D22' '2.108.70'
Set-TextValue 'E22' '  -1.77%  '
Set-TextValue 'D24' '7.460'
Set-TextValue 'E24' '  -1.54%  '
Set-TextValue 'D25' '162.11'
Set-TextValue 'E26' '  -0.07%  '
Set-TextValue 'D27' '0.1441'
Set-TextValue 'E27' '  -1.91%  '
Set-TextValue 'E28' '  +0.51%  '
Set-TextValue 'D29' '1.971'
Set-TextValue 'E29' '  +1.66%  '
Set-TextValue 'E30' '  +3.48%  '
Set-TextValue 'D31' '4.433'
Set-TextValue 'E31' '  -3.09%  '
Set-TextValue 'D32' '1.486'
Set-TextValue 'E32' '  -0.62%  '
Set-TextValue 'D33' '4.062'
Set-TextValue 'E33' '  +1.34%  '
Set-TextValue 'D34' '0.05217'
Set-TextValue 'E34' '  +0.75%  '
Set-TextValue 'E35' '  -1.33%  '
Set-TextValue 'D36' '0.7080'
Set-TextValue 'E36' '  +0.31%  '
Set-TextValue 'E37' '  -3.28%  '
Set-TextValue 'D38' '2.670'
Set-TextValue 'E38' '  +1.10%  '
Set-TextValue 'E39' '  -0.78%  '
Set-TextValue 'D40' '2.725'
Set-TextValue 'E40' '  +1.78%  '
Set-TextValue 'D41' '1.141.29'
Set-TextValue 'E41' '  +5.69%  '
Set-TextValue 'D42' '0.9183'
Set-TextValue 'E42' '  -1.82%  '
Set-TextValue 'D43' '5.976'
Set-TextValue 'E43' '  -0.44%  '
Set-TextValue 'D44' '0.4285'
Set-TextValue 'E44' '  -0.54%  '
Set-TextValue 'D45' '70.85'
Set-TextValue 'E45' '  +0.77%  '
Set-TextValue 'D46' '0.9998'
Set-TextValue 'E46' '  +0.05%  '
Set-TextValue 'D47' '102.67'
Set-TextValue 'E47' '  +0.33%  '
Set-TextValue 'E48' '  +1.13%  '
Set-TextValue 'D49' '2.005.96'
Set-TextValue 'E49' '  -0.46%  '
Set-TextValue 'D50' '9.181'
Set-TextValue 'E50' '  -0.52%  '
Set-TextValue 'D51' '6.975'
Set-TextValue 'E51' '  -1.38%  '
